$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old sheet had two stacked header rows (row 1 and row 2) above the data
# (rows 3..15). The new layout uses a single header row (row 1) followed
# directly by the (unchanged) data rows, which therefore shift up by one.
# Deleting the old row 2 accomplishes exactly that shift.
$ws.Rows.Item(2).Delete()

# Now rewrite row 1 with the new single-line header. A1:E1 use the plain
# default (unstyled) look, so reset any leftover formatting from the old
# header row first.
$left = $ws.Range("A1:E1")
$left.ClearFormats() | Out-Null
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Match the formatting used for the other header-style cells in the sheet
# (Arial 9pt, same font as the rest of the data cells).
$hdr = $ws.Range("F1:K1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 9

# Fix the selection so it matches the post-edit workbook state.
$ws.Range("A2:K2").Select() | Out-Null
